$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13, shifting old rows 13-24 down to 14-25
$ws.Rows(13).Insert()

# Row 13 previously had a label in A; the corrected layout has no A13 label
$ws.Range("A13").Clear()

# Establish correct formatting for the newly created B13/C13 cells
# by copying formats from a row that already has the right column styles
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Apply corrected text content
$ws.Range("B10").Value = "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."
$ws.Range("C10").Value = "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."

$ws.Range("B13").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C13").Value = "4808662 - Lucrécio Fábio dos Santos"

$ws.Range("B14").Value = "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"
$ws.Range("C14").Value = "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"

$ws.Range("B16").Value = ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."
$ws.Range("C16").Value = ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."

$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."

$ws.Range("B20").Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."

$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."

$ws.Range("B22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"
$ws.Range("C22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"
